$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.329.38'
$ws.Range('E2').Value = '  -0.95%  '
$ws.Range('D3').Value = '2.507.38'
$ws.Range('E3').Value = '  -2.18%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '305.53'
$ws.Range('E5').Value = '  +0.96%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '95.73'
$ws.Range('E6').Value = '  -1.57%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.583'
$ws.Range('E7').Value = '  +1.62%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.534'
$ws.Range('E9').Value = '  -2.28%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.41'
$ws.Range('E10').Value = '  +0.98%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0810'
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.65'
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('E13').Value = '  -3.56%  '
$ws.Range('D14').Value = '2.893.09'
$ws.Range('E14').Value = '  -2.35%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.32'
$ws.Range('E15').Value = '  +6.55%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.481.97'
$ws.Range('E16').Value = '  -1.30%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.856'
$ws.Range('E17').Value = '  -2.66%  '
$ws.Range('D18').Value = '42.347.64'
$ws.Range('E18').Value = '  -1.02%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.91'
$ws.Range('E19').Value = '  -1.99%  '
$ws.Range('D20').Value = '0.0₃0971'
$ws.Range('E20').Value = '  -1.46%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.41'
$ws.Range('E21').Value = '  -3.06%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '71.20'
$ws.Range('E22').Value = '  -0.48%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '252.77'
$ws.Range('E23').Value = '  -2.20%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.91'
$ws.Range('E24').Value = '  -1.41%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.03'
$ws.Range('E25').Value = '  -2.95%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '26.97'
$ws.Range('E26').Value = '  -3.64%  '
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.33'
$ws.Range('E28').Value = '  +10.84%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '10.14'
$ws.Range('E29').Value = '  +0.34%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '37.30'
$ws.Range('E30').Value = '  -4.59%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '5.92'
$ws.Range('E31').Value = '  -1.37%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '154.53'
$ws.Range('E32').Value = '  -1.34%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '19.23'
$ws.Range('E33').Value = '  +5.42%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.27'
$ws.Range('E34').Value = '  -1.16%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0784'
$ws.Range('E35').Value = '  -2.13%  '
$ws.Range('E36').Value = '  -3.63%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.60'
$ws.Range('E37').Value = '  -5.50%  '
$ws.Range('E38').Value = '  -0.78%  '
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '24.09'
$ws.Range('E39').Value = '  -1.87%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.119'
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('E42').Value = '  -0.42%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.999'
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('E44').Value = '  -2.57%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0299'
$ws.Range('E45').Value = '  -2.27%  '
$ws.Range('D46').Value = '2.030.82'
$ws.Range('E46').Value = '  -2.11%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '84.33'
$ws.Range('E47').Value = '  -4.99%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.00'
$ws.Range('E48').Value = '  -2.35%  '
$ws.Range('D49').Value = '2.749.44'
$ws.Range('E49').Value = '  -2.44%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '72.83'
$ws.Range('E50').Value = '  -5.94%  '
$ws.Range('E51').Value = '  -0.95%  '

